# Scheduled runner: refresh cached market-price columns (H:N) on Leve profit sheets.
# Values below come from a re-pull of current market data; only cells that actually
# moved are touched, and HQ-profit (N) / NQ-profit (M) cells that no longer apply are cleared
# to keep blank cells genuinely blank (matching the rest of the sheet).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 40 -- Leve Item ID (G40) = 5505
if ($ws.Range("G40").Value2 -eq 5505) {
    $ws.Range("H40").Value = 4512.75
    $ws.Range("I40").Value = 4166.6665
    $ws.Range("K40").Value = 4166.6665
    $ws.Range("M40").Value = -3991.6665
}

# Row 64 -- Leve Item ID (G64) = 5506
if ($ws.Range("G64").Value2 -eq 5506) {
    $ws.Range("H64").Value = 9006
    $ws.Range("J64").Value = 0
    $ws.Range("L64").Value = 0
    $ws.Range("N64").ClearContents()
}

# Row 67 -- Leve Item ID (G67) = 5506
if ($ws.Range("G67").Value2 -eq 5506) {
    $ws.Range("H67").Value = 9006
    $ws.Range("J67").Value = 0
    $ws.Range("L67").Value = 0
    $ws.Range("N67").ClearContents()
}

# Row 95 -- Leve Item ID (G95) = 18200
if ($ws.Range("G95").Value2 -eq 18200) {
    $ws.Range("H95").Value = 31623.5
    $ws.Range("J95").Value = 31623.5
    $ws.Range("L95").Value = 31623.5
    $ws.Range("N95").Value = -37115.5
}

# Row 137 -- Leve Item ID (G137) = 44013
if ($ws.Range("G137").Value2 -eq 44013) {
    $ws.Range("H137").Value = 1859
    $ws.Range("I137").Value = 800.4
    $ws.Range("J137").Value = 2615.1428
    $ws.Range("K137").Value = 2401.2
    $ws.Range("L137").Value = 7845.428400000001
    $ws.Range("M137").Value = 148.8000000000002
    $ws.Range("N137").Value = -12945.4284
}


$ws = $wb.Worksheets.Item("ARM")

# Row 101 -- Leve Item ID (G101) = 18518
if ($ws.Range("G101").Value2 -eq 18518) {
    $ws.Range("H101").Value = 0
    $ws.Range("J101").Value = 0
    $ws.Range("L101").Value = 0
    $ws.Range("N101").ClearContents()
}

# Row 102 -- Leve Item ID (G102) = 19945
if ($ws.Range("G102").Value2 -eq 19945) {
    $ws.Range("H102").Value = 42002300
    $ws.Range("I102").Value = 3335836.2
    $ws.Range("J102").Value = 100002000
    $ws.Range("K102").Value = 3335836.2
    $ws.Range("L102").Value = 100002000
    $ws.Range("M102").Value = -3334214.2
    $ws.Range("N102").Value = -100005244
}

# Row 132 -- Leve Item ID (G132) = 43997
if ($ws.Range("G132").Value2 -eq 43997) {
    $ws.Range("H132").Value = 3425
    $ws.Range("J132").Value = 2875
    $ws.Range("L132").Value = 8625
    $ws.Range("N132").Value = -13685
}


$ws = $wb.Worksheets.Item("BSM")

# Row 64 -- Leve Item ID (G64) = 14184
if ($ws.Range("G64").Value2 -eq 14184) {
    $ws.Range("H64").Value = 809.4
    $ws.Range("I64").Value = 786.3333
    $ws.Range("J64").Value = 844
    $ws.Range("K64").Value = 786.3333
    $ws.Range("L64").Value = 844
    $ws.Range("M64").Value = -561.3333
    $ws.Range("N64").Value = -1294
}

# Row 67 -- Leve Item ID (G67) = 14184
if ($ws.Range("G67").Value2 -eq 14184) {
    $ws.Range("H67").Value = 809.4
    $ws.Range("I67").Value = 786.3333
    $ws.Range("J67").Value = 844
    $ws.Range("K67").Value = 786.3333
    $ws.Range("L67").Value = 844
    $ws.Range("M67").Value = -6.333300000000008
    $ws.Range("N67").Value = -2404
}

# Row 94 -- Leve Item ID (G94) = 19939
if ($ws.Range("G94").Value2 -eq 19939) {
    $ws.Range("H94").Value = 221097.2
    $ws.Range("I94").Value = 367835.66
    $ws.Range("K94").Value = 367835.66
    $ws.Range("M94").Value = -367384.66
}

# Row 107 -- Leve Item ID (G107) = 27706
if ($ws.Range("G107").Value2 -eq 27706) {
    $ws.Range("H107").Value = 58567.57
    $ws.Range("I107").Value = 67510
    $ws.Range("K107").Value = 67510
    $ws.Range("M107").Value = -65590
}


$ws = $wb.Worksheets.Item("CRP")

# Row 28 -- Leve Item ID (G28) = 18348
if ($ws.Range("G28").Value2 -eq 18348) {
    $ws.Range("H28").Value = 21491
    $ws.Range("J28").Value = 21491
    $ws.Range("L28").Value = 21491
    $ws.Range("N28").Value = -21981
}

# Row 32 -- Leve Item ID (G32) = 2246
if ($ws.Range("G32").Value2 -eq 2246) {
    $ws.Range("H32").Value = 3499
    $ws.Range("I32").Value = 1873.5
    $ws.Range("K32").Value = 1873.5
    $ws.Range("M32").Value = -1557.5
}

# Row 58 -- Leve Item ID (G58) = 44021
if ($ws.Range("G58").Value2 -eq 44021) {
    $ws.Range("H58").Value = 1000
    $ws.Range("I58").Value = 1000
    $ws.Range("K58").Value = 1000
    $ws.Range("M58").Value = -797
}

# Row 107 -- Leve Item ID (G107) = 27689
if ($ws.Range("G107").Value2 -eq 27689) {
    $ws.Range("H107").Value = 506.22223
    $ws.Range("I107").Value = 276.16666
    $ws.Range("K107").Value = 276.16666
    $ws.Range("M107").Value = 1643.83334
}

# Row 132 -- Leve Item ID (G132) = 44019
if ($ws.Range("G132").Value2 -eq 44019) {
    $ws.Range("H132").Value = 1398.8334
    $ws.Range("J132").Value = 1622.5
    $ws.Range("L132").Value = 4867.5
    $ws.Range("N132").Value = -9927.5
}

# Row 136 -- Leve Item ID (G136) = 44021
if ($ws.Range("G136").Value2 -eq 44021) {
    $ws.Range("H136").Value = 1000
    $ws.Range("I136").Value = 1000
    $ws.Range("K136").Value = 3000
    $ws.Range("M136").Value = -450
}


$ws = $wb.Worksheets.Item("CUL")

# Row 16 -- Leve Item ID (G16) = 4641
if ($ws.Range("G16").Value2 -eq 4641) {
    $ws.Range("H16").Value = 33.25
    $ws.Range("I16").Value = 91
    $ws.Range("J16").Value = 14
    $ws.Range("K16").Value = 273
    $ws.Range("L16").Value = 42
    $ws.Range("M16").Value = -100
    $ws.Range("N16").Value = -388
}

# Row 63 -- Leve Item ID (G63) = 12866
if ($ws.Range("G63").Value2 -eq 12866) {
    $ws.Range("H63").Value = 1507.3334
    $ws.Range("I63").Value = 1011
    $ws.Range("K63").Value = 3033
    $ws.Range("M63").Value = -2284
}

# Row 66 -- Leve Item ID (G66) = 12866
if ($ws.Range("G66").Value2 -eq 12866) {
    $ws.Range("H66").Value = 1507.3334
    $ws.Range("I66").Value = 1011
    $ws.Range("K66").Value = 9099
    $ws.Range("M66").Value = -5355
}

# Row 68 -- Leve Item ID (G68) = 12895
if ($ws.Range("G68").Value2 -eq 12895) {
    $ws.Range("H68").Value = 2921.3076
    $ws.Range("I68").Value = 2997.3333
    $ws.Range("K68").Value = 8991.999899999999
    $ws.Range("M68").Value = -8180.999899999999
}

# Row 71 -- Leve Item ID (G71) = 12895
if ($ws.Range("G71").Value2 -eq 12895) {
    $ws.Range("H71").Value = 2921.3076
    $ws.Range("I71").Value = 2997.3333
    $ws.Range("K71").Value = 26975.9997
    $ws.Range("M71").Value = -22919.9997
}

# Row 108 -- Leve Item ID (G108) = 27853
if ($ws.Range("G108").Value2 -eq 27853) {
    $ws.Range("H108").Value = 3648.1667
    $ws.Range("I108").Value = 377.8
    $ws.Range("K108").Value = 1133.4
    $ws.Range("M108").Value = 1746.6
}

# Row 117 -- Leve Item ID (G117) = 27870
if ($ws.Range("G117").Value2 -eq 27870) {
    $ws.Range("H117").Value = 4156.3335
    $ws.Range("I117").Value = 736.2
    $ws.Range("J117").Value = 6599.2856
    $ws.Range("K117").Value = 2208.6
    $ws.Range("L117").Value = 19797.8568
    $ws.Range("M117").Value = 1233.4
    $ws.Range("N117").Value = -26681.8568
}

# Row 121 -- Leve Item ID (G121) = 27878
if ($ws.Range("G121").Value2 -eq 27878) {
    $ws.Range("H121").Value = 637.1111
    $ws.Range("I121").Value = 395.6
    $ws.Range("J121").Value = 939
    $ws.Range("K121").Value = 1186.8
    $ws.Range("L121").Value = 2817
    $ws.Range("M121").Value = 123.1999999999998
    $ws.Range("N121").Value = -5437
}


$ws = $wb.Worksheets.Item("GSM")

# Row 70 -- Leve Item ID (G70) = 14146
if ($ws.Range("G70").Value2 -eq 14146) {
    $ws.Range("H70").Value = 500000000
    $ws.Range("I70").Value = 0
    $ws.Range("K70").Value = 0
    $ws.Range("M70").ClearContents()
}

# Row 73 -- Leve Item ID (G73) = 14146
if ($ws.Range("G73").Value2 -eq 14146) {
    $ws.Range("H73").Value = 500000000
    $ws.Range("I73").Value = 0
    $ws.Range("K73").Value = 0
    $ws.Range("M73").ClearContents()
}

# Row 132 -- Leve Item ID (G132) = 44008
if ($ws.Range("G132").Value2 -eq 44008) {
    $ws.Range("H132").Value = 2415.25
    $ws.Range("I132").Value = 2331.7144
    $ws.Range("K132").Value = 6995.1432
    $ws.Range("M132").Value = -4465.1432
}


$ws = $wb.Worksheets.Item("LTW")

# Row 132 -- Leve Item ID (G132) = 44058
if ($ws.Range("G132").Value2 -eq 44058) {
    $ws.Range("H132").Value = 3490
    $ws.Range("I132").Value = 3490
    $ws.Range("J132").Value = 0
    $ws.Range("K132").Value = 10470
    $ws.Range("L132").Value = 0
    $ws.Range("M132").Value = -7940
    $ws.Range("N132").ClearContents()
}

# Row 136 -- Leve Item ID (G136) = 44060
if ($ws.Range("G136").Value2 -eq 44060) {
    $ws.Range("H136").Value = 3500
    $ws.Range("I136").Value = 3500
    $ws.Range("K136").Value = 10500
    $ws.Range("M136").Value = -7950
}


$ws = $wb.Worksheets.Item("WVR")

# Row 96 -- Leve Item ID (G96) = 19977
if ($ws.Range("G96").Value2 -eq 19977) {
    $ws.Range("H96").Value = 3112.875
    $ws.Range("I96").Value = 2833.8333
    $ws.Range("K96").Value = 2833.8333
    $ws.Range("M96").Value = -1460.8333
}

# Row 132 -- Leve Item ID (G132) = 44029
if ($ws.Range("G132").Value2 -eq 44029) {
    $ws.Range("H132").Value = 1728
    $ws.Range("I132").Value = 1507.76
    $ws.Range("J132").Value = 3563.3333
    $ws.Range("K132").Value = 4523.28
    $ws.Range("L132").Value = 10689.9999
    $ws.Range("M132").Value = -1993.28
    $ws.Range("N132").Value = -15749.9999
}

# Row 136 -- Leve Item ID (G136) = 44031
if ($ws.Range("G136").Value2 -eq 44031) {
    $ws.Range("H136").Value = 1992.7693
    $ws.Range("I136").Value = 1990.6364
    $ws.Range("K136").Value = 5971.9092
    $ws.Range("M136").Value = -3421.9092
}

